# Updates leads data for Istanbul pharmacies (rows 4-11), per commit:
# "✅ Veri Güncellendi: Istanbul"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Private-use-area glyph characters used as visual markers before the
# phone number / address text (preserved from the original cells).
$phoneIcon = [char]0xE0B0
$addrIcon  = [char]0xE0C8

# Row 4: Bahariye Cadde Eczanesi
$ws.Range("A4").Value = "Bahariye Cadde Eczanesi"
$ws.Range("B4").Value = "$phoneIcon`n+90 216 550 53 02"
$ws.Range("C4").Value = "$addrIcon`nOsmanağa, Gen. Asım Gündüz Caddesi No:17/C, 34734 Kadıköy/İstanbul, Türkiye"

# Row 5: İstanbul Eczanesi
$ws.Range("A5").Value = "İstanbul Eczanesi"
$ws.Range("B5").Value = "$phoneIcon`n+90 212 621 92 27"
$ws.Range("C5").Value = "$addrIcon`nHırka-i Şerif, Kocasinan Cd. no: 90, 34091 Fatih/İstanbul, Türkiye"

# Row 6: Eczane Yıldız
$ws.Range("A6").Value = "Eczane Yıldız"
$ws.Range("B6").Value = "$phoneIcon`n+90 543 842 68 58"
$ws.Range("C6").Value = "$addrIcon`nMecidiye Mahallesi, Ortaköy, Muvakkit Sk. No:10/A, 34347 Beşiktaş/İstanbul, Türkiye"

# Row 7: Colpan Pharmacy
$ws.Range("A7").Value = "Colpan Pharmacy"
$ws.Range("B7").Value = "$phoneIcon`n+90 212 523 56 63"
$ws.Range("C7").Value = "$addrIcon`nAkşemsettin, 34091 Fatih/İstanbul, Türkiye"

# Row 8 (Sultanahmet Eczanesi) is unchanged.

# Row 9: Istanbul Airport Pharmacy
$ws.Range("A9").Value = "Istanbul Airport Pharmacy"
$ws.Range("B9").Value = "$phoneIcon`n+90 212 830 35 57"
$ws.Range("C9").Value = "$addrIcon`nTayakadın Mah. İstanbul Yeni Havalimanı Giden Yolcu Katı NO: 7-F-0401, 34277 Arnavutköy/İstanbul, Türkiye"

# Row 10: İstanbul Havalimanı Eczane Melike Sultan
$ws.Range("A10").Value = "İstanbul Havalimanı Eczane Melike Sultan"
$ws.Range("B10").Value = "$phoneIcon`n+90 530 283 10 10"
$ws.Range("C10").Value = "$addrIcon`nİstanbul AIRPORT Havalimanı Dış Hatlar Gelen Katı No:820540, İmrahor, 34283 Arnavutköy/İstanbul, Türkiye"

# Row 11: Enes Eczanesi (was Eczapaketim / "Yok" / "Yok")
$ws.Range("A11").Value = "Enes Eczanesi"
$ws.Range("B11").Value = "$phoneIcon`n+90 530 283 10 10"
$ws.Range("C11").Value = "$addrIcon`nİstanbul AIRPORT Havalimanı Dış Hatlar Gelen Katı No:820540, İmrahor, 34283 Arnavutköy/İstanbul, Türkiye"
